# Applies the data rotation among rows 18, 19, 21, 22 (and a B-column-only
# update on row 20) on the "Artfynd" worksheet, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 18
$ws.Range("A18").Value = 112182654
$ws.Range("B18").Value = 77053
$ws.Range("E18").Value = 6437
$ws.Range("F18").Value = "Blanksvart spiklav"
$ws.Range("G18").Value = "Calicium denigratum"
$ws.Range("H18").Value = "(Vain.) Tibell"
$ws.Range("Q18").Value = 364914
$ws.Range("R18").Value = 6872133

# Row 19
$ws.Range("A19").Value = 112182890
$ws.Range("B19").Value = 96652
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 219790
$ws.Range("F19").Value = "Fläcknycklar"
$ws.Range("G19").Value = "Dactylorhiza maculata"
$ws.Range("H19").Value = "(L.) Soó"
$ws.Range("Q19").Value = 364947
$ws.Range("R19").Value = 6872308

# Row 20
$ws.Range("B20").Value = 78714

# Row 21
$ws.Range("A21").Value = 112182494
$ws.Range("B21").Value = 77402
$ws.Range("E21").Value = 6446
$ws.Range("F21").Value = "Kolflarnlav"
$ws.Range("G21").Value = "Carbonicola anthracophila"
$ws.Range("H21").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q21").Value = 364938
$ws.Range("R21").Value = 6872236

# Row 22
$ws.Range("A22").Value = 112182349
$ws.Range("B22").Value = 77650
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("Q22").Value = 364898
$ws.Range("R22").Value = 6872201
